$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New OHLCV rows (row,A,B,C,D,E,F) to append starting at row 1085
$newRowsData = @(
"1085,45534.5,2512.43,2552.17,2434.96,2448.75,103983.6261",
"1086,45534.66666666666,2448.76,2542.19,2431.14,2506.94,130701.1151",
"1087,45534.83333333334,2506.94,2530.66,2506.35,2526,19268.6521",
"1088,45535,2525.99,2533.95,2518.86,2523.52,14166.9397",
"1089,45535.16666666666,2523.52,2529.86,2518.17,2521.01,14956.1787",
"1090,45535.33333333334,2521.01,2529.81,2511.47,2525.5,14773.6726",
"1091,45535.5,2525.49,2527.8,2513.8,2519.19,15368.1774",
"1092,45535.66666666666,2519.18,2519.8,2491.92,2499.8,22057.3239",
"1093,45535.83333333334,2499.79,2512.64,2493.68,2510.89,7778.3142",
"1094,45536,2513,2516.28,2476.44,2488.74,22040.3526",
"1095,45536.16666666666,2488.75,2491.18,2458.84,2486.82,26820.4775",
"1096,45536.33333333334,2486.82,2492.88,2459.02,2477.74,22368.1709",
"1097,45536.5,2477.73,2489.31,2438,2480.83,47978.3656",
"1098,45536.66666666666,2480.83,2487.01,2466.54,2476.14,7689.2003",
"1099,45536.83333333334,2501.71,2516,2473.6,2482.94,11700.672",
"1100,45537,2425.71,2455.26,2423.52,2450.52,35331.6569",
"1101,45537.16666666666,2450.47,2466.2,2433.18,2441.14,37943.8673",
"1102,45537.33333333334,2441.14,2539.53,2438.39,2517.33,84869.04240000001",
"1103,45537.5,2515.81,2524.67,2510.59,2520.02,7583.8207",
"1104,45537.66666666666,2513.2,2530.92,2505.78,2524.25,29916.1001",
"1105,45537.83333333334,2524.26,2564.83,2522.22,2538.01,31904.11",
"1106,45538,2538,2553.6,2514.73,2516,27730.1585",
"1107,45538.16666666666,2516,2525.2,2512.46,2519.55,19383.1812",
"1108,45538.33333333334,2519.55,2520.54,2494.24,2509.59,20704.4681",
"1109,45538.5,2509.59,2519.12,2435.79,2444.01,81518.50019999999",
"1110,45538.66666666666,2444.01,2469.87,2435.83,2447.6,41528.1331",
"1111,45538.83333333334,2447.6,2470.57,2411.12,2425.29,33928.6888",
"1112,45539,2425.28,2449.22,2306.65,2373,144366.9266",
"1113,45539.16666666666,2372.99,2414.71,2363.56,2401.52,45719.3498",
"1114,45539.33333333334,2401.52,2410.83,2387.15,2398.77,36620.6349",
"1115,45539.5,2398.77,2458.57,2393.16,2455.4,91240.01700000001",
"1116,45539.66666666666,2455.39,2490,2439.2,2447.62,46355.3428",
"1117,45539.83333333334,2452.15,2469.55,2443.32,2450.71,18010.6716",
"1118,45540,2450.71,2466,2404.28,2405.55,51182.1435",
"1119,45540.16666666666,2405.56,2418.47,2376.72,2412.4,42138.4727",
"1120,45540.33333333334,2412.41,2416.35,2383.16,2389.54,28448.2371",
"1121,45540.5,2389.53,2415.15,2355.8,2371.43,85446.0382",
"1122,45540.66666666666,2371.43,2395.99,2348.04,2363.97,43189.7973",
"1123,45540.83333333334,2363.97,2378.33,2350,2368.81,28144.4081",
"1124,45541,2368.81,2408.4,2357.6,2392.63,32980.5199",
"1125,45541.16666666666,2392.63,2397.35,2321.21,2344.45,57143.7515",
"1126,45541.33333333334,2344.44,2377.4,2338.57,2372.8,36067.1436"
)

$startRow = 1085
$lastRow = $startRow + $newRowsData.Count - 1

# Build a 2D array for fast bulk write of columns A:F
$values = New-Object 'object[,]' $newRowsData.Count, 6

for ($i = 0; $i -lt $newRowsData.Count; $i++) {
    $parts = $newRowsData[$i].Split(",")
    $values[$i, 0] = [double]$parts[1]
    $values[$i, 1] = [double]$parts[2]
    $values[$i, 2] = [double]$parts[3]
    $values[$i, 3] = [double]$parts[4]
    $values[$i, 4] = [double]$parts[5]
    $values[$i, 5] = [double]$parts[6]
}

$targetRange = $ws.Range("A" + $startRow + ":F" + $lastRow)
$targetRange.Value = $values

# Match the style used by the other date cells in column A (bold, centered,
# bordered, custom date number format "YYYY-MM-DD HH:MM:SS") by copying the
# format from the row immediately above the newly appended block. Using
# Copy/PasteSpecial (instead of rebuilding the format property-by-property)
# reuses the existing style record instead of minting new ones.
$styleSource = $ws.Range("A" + ($startRow - 1))
$styleTarget = $ws.Range("A" + $startRow + ":A" + $lastRow)
$styleSource.Copy()
$styleTarget.PasteSpecial(-4122)

Write-Output ("Appended rows {0} to {1}" -f $startRow, $lastRow)
